# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.666.19'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '1.884.80'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'249.51"
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = "'0.4763"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = "'0.2942"
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").Value = "'0.06548"
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = "'22.06"
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = "'0.7411"
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = "'96.88"
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '1.882.74'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = "'5.250"
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").Value = "'275.64"
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '30.641.80'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = "'13.20"
$ws.Range("E18").Value = '  -3.07%  '
$ws.Range("D19").Value = "'0.000007560"
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").Value = '2.131.07'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = "'5.349"
$ws.Range("E22").Value = '  +1.72%  '
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = "'6.245"
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").Value = "'9.250"
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").Value = "'164.31"
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'1.918"
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("D30").Value = "'0.09743"
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = "'4.296"
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").Value = "'4.175"
$ws.Range("D34").Value = "'0.04916"
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("D35").Value = "'1.128"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = "'0.7012"
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = "'2.726"
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("D38").Value = "'0.01914"
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").Value = "'6.325"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").Value = "'75.72"
$ws.Range("E41").Value = '  +6.17%  '
$ws.Range("D42").Value = "'2.032"
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("D43").Value = "'0.4261"
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("D44").Value = "'0.8434"
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = "'102.76"
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = "'9.439"
$ws.Range("E47").Value = '  +1.92%  '
$ws.Range("D48").Value = "'7.069"
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").Value = "'35.68"
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").Value = "'918.57"
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("E51").Value = '  +2.13%  '
